# Move the two names "8fcbyk1m" and "y1puz53f" from the "Sheet1" (unused
# names) list to the "used" tracking sheet, recording the new image
# filenames + usage timestamps, per the commit:
#   "Update works images 2026-01-21 15:48:30"

$wb = $excel.ActiveWorkbook

$namesSheet = $wb.Worksheets.Item("Sheet1")
$usedSheet  = $wb.Worksheets.Item("used")

# --- Sheet1: remove the two consumed names (rows 2 and 3), shifting
# the remaining names up so the list stays contiguous. ---
$namesSheet.Range("A2:A3").EntireRow.Delete()

# --- used: append the two names that were just consumed, with their
# source filename and the timestamp they were used. ---
$lastRow = $usedSheet.UsedRange.Rows.Count

$usedSheet.Range("A" + ($lastRow + 1)).Value = "8fcbyk1m"
$usedSheet.Range("B" + ($lastRow + 1)).Value = "ChatGPT Image 2026年1月21日 15_47_37.png"
$usedSheet.Range("C" + ($lastRow + 1)).Value = "2026-01-21 15:48:23"

$usedSheet.Range("A" + ($lastRow + 2)).Value = "y1puz53f"
$usedSheet.Range("B" + ($lastRow + 2)).Value = "ChatGPT Image 2026年1月21日 15_47_41.png"
$usedSheet.Range("C" + ($lastRow + 2)).Value = "2026-01-21 15:48:23"
